# Update New Orleans shard workbook:
#  1. hotel_info gains a new "State" column (inserted right after "Hotel_Name"),
#     populated with "Louisiana" for the existing data row.
#  2. The worksheet tab order is swapped so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# Insert a new column C (State) on hotel_info, shifting City/Zip/etc. right.
$hotelSheet.Range("C1:C2").EntireColumn.Insert()
$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"

# Reorder tabs: review_info first, hotel_info second.
$reviewSheet.Move($wb.Worksheets.Item(1))
